# "Generate Report for Handoff":
#   - Status text "Handed back: in sync with en-US" -> "Ready for handoff"
#     on the Overview sheet (zh-cn/de-de status columns) and on each
#     per-language sheet's Status column.
#   - The corresponding "Latest Handoff/Generate" timestamps are bumped
#     forward by the report-regeneration run.
#   - The Status columns are narrowed (they no longer need to fit the old,
#     longer "Handed back: in sync with en-US" text).
#     Note: ColumnWidth is stored/quantized by the engine on a fixed pixel
#     grid, so we pick the input that snaps to the closest achievable
#     stored width to the target.
$wb = $excel.ActiveWorkbook

# --- Sheet: Overview ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-29 17:06:44"
$wsOverview.Range("E1").ColumnWidth = 16.333333333333336
$wsOverview.Range("F1").ColumnWidth = 16.333333333333336

# --- Sheet: zh-cn ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-29 17:06:38"
$wsZhCn.Range("C1").ColumnWidth = 16.333333333333336

# --- Sheet: de-de ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-29 17:06:44"
$wsDeDe.Range("C1").ColumnWidth = 16.333333333333336
